$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix typo "tempratuur" -> "temperatuur" in the relevant cells
$ws.Range("B2").Value = "temperatuur sensor"
$ws.Range("C2").Value = "hoeveel de temperatuur sensor af mag wijken van de werkelijke temperatuur"
$ws.Range("B7").Value = "temperatuur weergeven"
$ws.Range("C7").Value = "snel en hoe vaak moet de temperatuur weergaven geüpdate worden"

# Update the active selection to B2
$ws.Activate()
$ws.Range("B2").Select()
